# "pages and logo up" - fill in the author's name on the title page and
# clean up the stray _GoBack bookmark that Word leaves behind after an edit.

$d = $word.ActiveDocument

# 1. Replace the "Student Name" placeholder on the title page with the
#    author's actual name, editing the paragraph range directly so the
#    run's existing formatting/attributes are preserved.
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*Student Name*") {
        $p.Range.Text = "Tina Bekkholt"
    }
}

# 2. Remove the leftover "_GoBack" bookmark (an artifact Word inserts at
#    the last edit position) that was cleaned up in this revision.
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}
